$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "is_locked_lbl" (D1) and "is_enabled_lbl" (E1) template columns,
# shifting order_by/rem left into D1/E1.
$ws.Range("D1:E1").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)

# Append a new "tenant_id_lbl" template column after rem (now in F1).
$ws.Range("F1").Value = '<%=comment.tenant_id_lbl%><%selectList.tenant_id = data.findAllTenant.map((item) => item.lbl)%><%_dataValidation_({ sqref: `${ _col }2:${ _col }${ _lastRow }`, formula1: `"${ selectList.tenant_id.join(",") }"` })%>'
